# calorimetry : scripts : tests : updated
#
# 1. Rename sheet "adj_r_squared" -> "metrics"
# 2. Rebuild that sheet's data as a small metrics table (Adj.R^2, NRMSE, SMAPE, RMSE)
# 3. Clear the stray extra rows on "input_enthalpies" (keep only the header row)

$wb = $excel.ActiveWorkbook

# --- 1 & 2: adj_r_squared -> metrics ---------------------------------------
$ws = $wb.Worksheets.Item("adj_r_squared")
$ws.Name = "metrics"

$ws.Range("A1").Value = "metrics"
$ws.Range("B1").Value = "value"

$ws.Range("A2").Value = "Adj.R^2"
$ws.Range("B2").Value = 0.994740757337271

$ws.Range("A3").Value = "NRMSE"
$ws.Range("B3").Value = 0.060744579194617

$ws.Range("A4").Value = "SMAPE"
$ws.Range("B4").Value = 0.516709138142194

$ws.Range("A5").Value = "RMSE"
$ws.Range("B5").Value = 0.00221468675054425

# --- 3: input_enthalpies loses its extra (blank / zero-dev) rows -----------
$ie = $wb.Worksheets.Item("input_enthalpies")
$ie.Rows.Item(2).EntireRow.Delete()
$ie.Rows.Item(2).EntireRow.Delete()

Write-Host "metrics sheet + input_enthalpies cleanup applied"
